# Applies the "Updated cryptos list" data refresh to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that looks numeric (e.g. "573.68", "64.511.64").
# Excel auto-converts such literals typed into a General-formatted cell into
# real numbers, which would change the stored cell type from text to number.
# Force the column to Text format first so values round-trip as strings,
# then restore the original (unformatted / style-less) "Normal" style so no
# stray number-format is left behind on the cells.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '64.511.64'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '3.139.37'
$ws.Range('E3').Value = '  -1.29%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '573.68'
$ws.Range('E5').Value = '  +0.24%  '
$ws.Range('D6').Value = '164.33'
$ws.Range('E6').Value = '  -3.19%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = '0.575'
$ws.Range('E8').Value = '  -5.80%  '
$ws.Range('D9').Value = '3.153.56'
$ws.Range('E9').Value = '  -1.12%  '
$ws.Range('E10').Value = '  -2.59%  '
$ws.Range('D11').Value = '6.64'
$ws.Range('E11').Value = '  -3.04%  '
$ws.Range('D12').Value = '0.385'
$ws.Range('E12').Value = '  -1.69%  '
$ws.Range('D13').Value = '3.689.98'
$ws.Range('E13').Value = '  -1.42%  '
$ws.Range('D14').Value = '0.127'
$ws.Range('E14').Value = '  -1.59%  '
$ws.Range('D15').Value = '64.505.21'
$ws.Range('E15').Value = '  -0.09%  '
$ws.Range('D16').Value = '25.00'
$ws.Range('E16').Value = '  -1.49%  '
$ws.Range('D17').Value = '3.147.76'
$ws.Range('E17').Value = '  -1.44%  '
$ws.Range('E18').Value = '  -2.44%  '
$ws.Range('D19').Value = '407.66'
$ws.Range('E19').Value = '  -3.19%  '
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').Value = '5.24'
$ws.Range('E20').Value = '  -2.46%  '
$ws.Range('B21').Value = 'Chainlink'
$ws.Range('C21').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D21').Value = '12.53'
$ws.Range('E21').Value = '  -3.86%  '
$ws.Range('D22').Value = '7.07'
$ws.Range('E22').Value = '  -1.63%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').Value = '68.91'
$ws.Range('E24').Value = '  -2.05%  '
$ws.Range('D25').Value = '0.484'
$ws.Range('E25').Value = '  -2.83%  '
$ws.Range('E26').Value = '  -5.42%  '
$ws.Range('E27').Value = '  -2.78%  '
$ws.Range('D28').Value = '8.94'
$ws.Range('E28').Value = '  +1.78%  '
$ws.Range('D29').Value = '0.997'
$ws.Range('E29').Value = '  +0.19%  '
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('D31').Value = '1.81'
$ws.Range('E31').Value = '  -1.69%  '
$ws.Range('D32').Value = '21.25'
$ws.Range('E32').Value = '  -2.67%  '
$ws.Range('D33').Value = '163.05'
$ws.Range('E33').Value = '  +3.62%  '
$ws.Range('E34').Value = '  -4.17%  '
$ws.Range('D35').Value = '6.29'
$ws.Range('E35').Value = '  -1.43%  '
$ws.Range('E36').Value = '  +0.25%  '
$ws.Range('D37').Value = '1.36'
$ws.Range('E37').Value = '  -0.40%  '
$ws.Range('E38').Value = '  -1.54%  '
$ws.Range('D39').Value = '2.635.60'
$ws.Range('E39').Value = '  -3.10%  '
$ws.Range('D40').Value = '23.75'
$ws.Range('E40').Value = '  -2.53%  '
$ws.Range('D41').Value = '4.10'
$ws.Range('E41').Value = '  -3.58%  '
$ws.Range('D42').Value = '38.21'
$ws.Range('E42').Value = '  -2.42%  '
$ws.Range('D43').Value = '0.692'
$ws.Range('E43').Value = '  -3.54%  '
$ws.Range('D44').Value = '0.0614'
$ws.Range('E44').Value = '  -1.59%  '
$ws.Range('D45').Value = '5.34'
$ws.Range('E45').Value = '  -4.22%  '
$ws.Range('D46').Value = '290.22'
$ws.Range('E46').Value = '  -0.74%  '
$ws.Range('D47').Value = '21.38'
$ws.Range('E47').Value = '  -0.50%  '
$ws.Range('D48').Value = '0.0254'
$ws.Range('E48').Value = '  -3.38%  '
$ws.Range('E49').Value = '  -0.15%  '
$ws.Range('D50').Value = '0.0974'
$ws.Range('E50').Value = '  -2.10%  '
$ws.Range('D51').Value = '10.49'
$ws.Range('E51').Value = '  +0.43%  '

$priceRange.Style = "Normal"
